$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add headers for new columns I and J, matching the existing header style (H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for I2:I57 and J2:J57
$iValues = @(11,5,7,7,5,5,1,7,6,9,1,1,5,9,9,5,8,7,6,7,9,9,8,7,8,8,9,4,5,8,7,9,5,9,2,7,9,8,7,9,2,9,8,7,9,5,5,6,8,8,6,2,6,5,6,6)
$jValues = @(12,5,8,7,5,6,2,8,6,9,3,3,6,9,9,6,9,8,7,8,9,9,8,8,8,8,9,5,7,8,7,9,5,9,3,7,9,8,8,9,3,9,8,8,9,6,6,6,8,8,6,2,6,5,6,6)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}
